$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Market Cap (column C) values for rows 2-26
$ws.Range("C2").Value = 721202497449.1852
$ws.Range("C3").Value = 244685368352.3357
$ws.Range("C4").Value = 37832035877.15541
$ws.Range("C5").Value = 33930139813.33422
$ws.Range("C6").Value = 26502594024.89571
$ws.Range("C7").Value = 13932760853.96851
$ws.Range("C8").Value = 11143606318.32923
$ws.Range("C9").Value = 9172445970.432299

# Rows 10 and 11 swap Name/Symbol (Polygon now before Avalanche)
$ws.Range("A10").Value = "Polygon"
$ws.Range("B10").Value = "MATIC-USD"
$ws.Range("C10").Value = 8288439414.011463

$ws.Range("A11").Value = "Avalanche"
$ws.Range("B11").Value = "AVAX-USD"
$ws.Range("C11").Value = 8286989881.493361

$ws.Range("C12").Value = 8057833605.645318
$ws.Range("C13").Value = 8005791334.873327
$ws.Range("C14").Value = 7089270478.055191
$ws.Range("C15").Value = 6024677236.291335
$ws.Range("C16").Value = 5363446399.970418
$ws.Range("C17").Value = 5206001255.311966
$ws.Range("C18").Value = 4595203878.57365
$ws.Range("C19").Value = 3719839375.908229
$ws.Range("C20").Value = 3706033919.912179
$ws.Range("C21").Value = 3521926218.888741
$ws.Range("C22").Value = 3383257432.909583
$ws.Range("C23").Value = 3017812634.355353
$ws.Range("C24").Value = 2951921492.128325
$ws.Range("C25").Value = 2881111534.018945
$ws.Range("C26").Value = 2525262868.699377
